# Update "想去人数" (interested-count) figures in the 展览 and 全部类型 sheets
# to reflect newly generated data (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 6618
$wsExhibit.Range("F4").Value = 106
$wsExhibit.Range("F5").Value = 152
$wsExhibit.Range("F8").Value = 588

# Sheet "全部类型" (All types) - aggregated view containing the same rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6618
$wsAll.Range("F5").Value = 106
$wsAll.Range("F6").Value = 152
$wsAll.Range("F10").Value = 588
